$d = $word.ActiveDocument

# The text "<id>p024r_1</id>" is split across three runs in the source:
#   1) "<id>"      (Courier New, color 7f6000, sz 18)
#   2) "p024r_1"   (plain, color 000000)
#   3) "</id>"     (Courier New, color 7f6000, sz 18)
# Locate it and collapse it into a single run, keeping run 1's formatting,
# resulting in one run whose text is "<id>p024r_1</id>".

$needle = "<id>p024r_1</id>"

$r = $d.Content
$found = $r.Find.Execute($needle, $false, $false, $false, $false, $false,
                          $true, 1, $false, "", 0)

$start = $r.Start
$end = $r.End
$firstRunEnd = $start + 4   # length of "<id>"

$restRange = $d.Range($firstRunEnd, $end)
$restText = $restRange.Text
$restRange.Delete()

$firstRunRange = $d.Range($start, $firstRunEnd)
$firstRunRange.InsertAfter($restText)
